$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("StartSceneConfig")

# Update column H values in rows 6..36: 10002..10032 -> 30001..30031
for ($row = 6; $row -le 36; $row++) {
    $newValue = 30001 + ($row - 6)
    $ws.Cells.Item($row, 8).Value = $newValue
}

# Update the selected cell in the sheet view from B8 to I16
$ws.Range("I16").Select()
